$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 26 and row 27 for columns A, B, E, F, G, H, Q, R
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell26 = $ws.Range($col + "26")
    $cell27 = $ws.Range($col + "27")

    $val26 = $cell26.Value2
    $val27 = $cell27.Value2

    $cell26.Value2 = $val27
    $cell27.Value2 = $val26
}
